$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 197
$wsExpo.Range("F4").Value = 2419
$wsExpo.Range("F5").Value = 30

# Update "全部类型" (All types) sheet - mirrors the same rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 197
$wsAll.Range("F6").Value = 2419
$wsAll.Range("F7").Value = 30
